# Tenta Avancerad Programmering.docx - apply the "Sakerhet" paragraph edit:
#   - append " Authorization för microservices." to the bold "Säkerhet: " paragraph,
#     splitting the run so that only "Säkerhet: " stays bold and the appended text
#     is regular weight, with proofErr spell-check markers wrapping "Authorization"
#   - the trailing hidden "_GoBack" bookmark (Word's "last edit location" marker)
#     moves from the end of the previous paragraph to the end of this paragraph,
#     since that's where the new text was typed.

$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the end of the previous paragraph
# ("... detta, samt läsa mer information."). Remove it first so that re-adding
# it below (after the new text) gets bookmark id 0 again, matching the target.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the "Säkerhet: " paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Säkerhet: ") {
        $target = $cand
        break
    }
}

$r = $target.Range

# Rebuild the paragraph's contents in place via WordprocessingML, preserving the
# paragraph's own identity attributes, dropping the paragraph-mark bold (pPr/rPr)
# now that real (non-mark) content follows it, splitting "Säkerhet:" / " " into
# their own bold runs, and appending the new, non-bold text with spell-check
# proofErr markers around "Authorization", followed by the relocated bookmark.
$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D56C19" w:rsidRPr="00EE4854" w:rsidRDefault="00D56C19" w:rsidP="00B7158B"><w:r><w:rPr><w:b/></w:rPr><w:t>Säkerhet:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Authorization</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> för microservices.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($xmlFrag)

Write-Host "Updated paragraph text:" $target.Range.Text
